$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 405, shifting rows 405:460 down to 406:461.
$ws.Rows.Item(405).Insert()

# Populate the newly inserted row 405 with the new record's data.
$ws.Cells.Item(405, 1).Value = 5
$ws.Cells.Item(405, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(405, 3).Value = "Maule"
$ws.Cells.Item(405, 4).Value = 45124
$ws.Cells.Item(405, 5).Value = 7
$ws.Cells.Item(405, 6).Value = 100112008
$ws.Cells.Item(405, 7).Value = "Coliflor"
$ws.Cells.Item(405, 8).Value = "Sin especificar"
$ws.Cells.Item(405, 9).Value = "Primera"
$ws.Cells.Item(405, 10).Value = 4000
$ws.Cells.Item(405, 11).Value = 600
$ws.Cells.Item(405, 12).Value = 600
$ws.Cells.Item(405, 13).Value = 600
$ws.Cells.Item(405, 14).Value = "$/unidad"
$ws.Cells.Item(405, 15).Value = "Región del Maule"
$ws.Cells.Item(405, 16).Value = 600
$ws.Cells.Item(405, 17).Value = 1
$ws.Cells.Item(405, 18).Value = "Hortaliza"
